# RatingScore.xlsx edit script
# - Removes the "Added by Steven Zhang, needs to confirm with Daphne" placeholder note
# - Inserts a proper "BBB *-" / Fitch row (with a "Same as "BBB"" note) right after the
#   existing "BBB" Fitch row, instead of leaving it as a stray highlighted row at the
#   bottom of the table
# - Clears the old bottom "BBB *-" row down to a blank trailing spacer row
# - Replaces the old ad-hoc comment on the Moody's "B1u" row with a proper explanation
#   of Moody's "u"/"e" rating identifiers, and clears the stray note on the "Caa1u" row
# - Removes the yellow highlight fill that was used to flag these rows for review
# - Updates the view position/selection to where the edit was made

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Insert a new row before row 73 (the "Fitch BBB-" row) ------------------------
# This pushes the old rows 73..102 down to 74..103, which conveniently carries the old
# trailing "BBB *-" row (formerly row 102) down to row 103 automatically.
$ws.Rows.Item(73).Insert()

# --- 2) Populate the newly inserted row 73 with the Fitch "BBB *-" entry -------------
$ws.Cells.Item(73, 1).Value2 = "Fitch"
$ws.Cells.Item(73, 2).Value2 = "BBB *-"
$ws.Cells.Item(73, 3).Value2 = 13
$ws.Cells.Item(73, 6).Value2 = 'Same as "BBB"'

# --- 3) Clean up the old trailing row (now row 103): drop the Fitch/"BBB *-" data and
#        the old placeholder comment, leaving only the blank styled spacer cells -------
$ws.Range("A103:C103").ClearContents()
$ws.Range("A103:C103").ClearFormats()
$ws.Cells.Item(103, 6).ClearContents()

# --- 4) Replace the ad-hoc reviewer notes on the Moody's "u" rows -------------------
# Row 45: Moody's "B1u" -- swap the placeholder note for a real explanation
$ws.Cells.Item(45, 6).Value2 = "In recent years, Moody" + [char]8217 + "s added " + [char]8220 + "u" + [char]8221 + "  or " + [char]8220 + "e" + [char]8221 + " identifier to express additional meanings on their ratings. These identifiers do not change the associated rating score actually. "

# Row 49: Moody's "Caa1u" -- just clear the placeholder note
$ws.Cells.Item(49, 6).ClearContents()

# --- 5) Remove the yellow review-highlight fill from all flagged rows ---------------
$ws.Range("A45:G45").Interior.Pattern = -4142
$ws.Range("A49:G49").Interior.Pattern = -4142
$ws.Range("A73:C73").Interior.Pattern = -4142
$ws.Range("D103:I103").Interior.Pattern = -4142

# --- 6) Update the view: scroll position + active selection -------------------------
$ws.Range("H73").Select()
